$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '66.297.80'
$ws.Range('E2').Value = '  +0.17%  '
Set-TextValue 'D3' '3.590.83'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  +0.20%  '
Set-TextValue 'D5' '605.65'
$ws.Range('E5').Value = '  -0.14%  '
Set-TextValue 'D6' '148.03'
$ws.Range('E6').Value = '  +2.48%  '
Set-TextValue 'D7' '3.592.49'
$ws.Range('E7').Value = '  +0.59%  '
$ws.Range('E8').Value = '  -0.05%  '
Set-TextValue 'D9' '0.486'
$ws.Range('E9').Value = '  -0.77%  '
Set-TextValue 'D10' '0.135'
$ws.Range('E10').Value = '  -0.67%  '
Set-TextValue 'D11' '7.86'
$ws.Range('E11').Value = '  +0.58%  '
Set-TextValue 'D12' '0.414'
$ws.Range('E12').Value = '  +0.38%  '
Set-TextValue 'D13' '4.206.88'
$ws.Range('E13').Value = '  +0.84%  '
Set-TextValue 'D14' '0.0000205'
$ws.Range('E14').Value = '  -1.33%  '
Set-TextValue 'D15' '29.55'
$ws.Range('E15').Value = '  -2.00%  '
Set-TextValue 'D16' '3.595.30'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('E17').Value = '  +1.93%  '
Set-TextValue 'D18' '66.426.43'
$ws.Range('E18').Value = '  +0.26%  '
Set-TextValue 'D19' '11.09'
$ws.Range('E19').Value = '  -3.00%  '
Set-TextValue 'D20' '6.31'
$ws.Range('E20').Value = '  +1.40%  '
Set-TextValue 'D21' '14.89'
$ws.Range('E21').Value = '  +1.02%  '
Set-TextValue 'D22' '422.17'
$ws.Range('E22').Value = '  -2.22%  '
Set-TextValue 'D23' '0.611'
$ws.Range('E23').Value = '  +0.12%  '
Set-TextValue 'D24' '78.44'
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('E25').Value = '  -0.10%  '
Set-TextValue 'D26' '0.0000119'
$ws.Range('E26').Value = '  +0.87%  '
Set-TextValue 'D27' '8.23'
$ws.Range('E27').Value = '  +3.99%  '
Set-TextValue 'D28' '9.34'
$ws.Range('E28').Value = '  +2.26%  '
Set-TextValue 'D29' '2.48'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('E30').Value = '  +0.30%  '
Set-TextValue 'D31' '3.592.15'
$ws.Range('E31').Value = '  +0.74%  '
Set-TextValue 'D32' '0.156'
$ws.Range('E32').Value = '  +2.59%  '
Set-TextValue 'D33' '1.43'
$ws.Range('E33').Value = '  -1.45%  '
Set-TextValue 'D34' '25.03'
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('E35').Value = '  +0.01%  '
Set-TextValue 'D36' '7.73'
$ws.Range('E36').Value = '  -1.51%  '
Set-TextValue 'D37' '5.59'
$ws.Range('E37').Value = '  +0.27%  '
Set-TextValue 'D38' '1.67'
$ws.Range('E38').Value = '  -3.23%  '
Set-TextValue 'D39' '174.65'
$ws.Range('E39').Value = '  +0.76%  '
Set-TextValue 'D40' '0.0844'
$ws.Range('E40').Value = '  -0.10%  '
Set-TextValue 'D41' '5.17'
$ws.Range('E41').Value = '  -0.72%  '
Set-TextValue 'D42' '0.888'
$ws.Range('E42').Value = '  -0.94%  '
Set-TextValue 'D43' '45.85'
$ws.Range('E43').Value = '  -0.36%  '
Set-TextValue 'D44' '1.84'
$ws.Range('E44').Value = '  -6.18%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D45' '1.00'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D46' '2.51'
$ws.Range('E46').Value = '  +4.17%  '
Set-TextValue 'D47' '24.15'
$ws.Range('E47').Value = '  -3.87%  '
$ws.Range('B48').Value = 'Cosmos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D48' '7.13'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D49' '23.44'
$ws.Range('E49').Value = '  +1.23%  '
Set-TextValue 'D50' '1.13'
$ws.Range('E50').Value = '  -6.25%  '
Set-TextValue 'D51' '0.957'
$ws.Range('E51').Value = '  +2.44%  '
